$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data area (A2:C11) before rewriting with the new layout
$ws.Range("A2:C11").Clear() | Out-Null

# Seed the two new shared strings in the same order the original workbook
# author created them (Left Encoder A/B, then Right Encoder A/B) so the
# resulting sharedStrings table matches exactly, even though "Right Encoder"
# rows are laid out above "Left Encoder" rows in the final sheet.
$ws.Cells.Item(11, 1).Value = "Left Encoder A"
$ws.Cells.Item(12, 1).Value = "Left Encoder B"

# New pin configuration (rows 2-15)
$data = @(
    @("L Motor Speed",   7,  $null),
    @("L Motor Reverse", 8,  $null),
    @("L Motor Forward", 9,  $null),
    @("R Motor Speed",   10, $null),
    @("R Motor Forward", 11, $null),
    @("R Motor Reverse", 12, $null),
    @("STBY",            24, $null),
    @("Right Encoder A", 27, $null),
    @("Right Encoder B", 28, $null),
    @("Left Encoder A",  37, $null),
    @("Left Encoder B",  38, $null),
    @("IR Right",        32, "A13"),
    @("IR Left",         33, "A14"),
    @("IR Middle",       34, "A15")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    if ($entry[2] -ne $null) {
        $ws.Cells.Item($row, 3).Value = $entry[2]
    }
    $row++
}

# Column B needs an explicit width now that it holds wider numeric pin values
$ws.Columns.Item(2).ColumnWidth = 7.28515625

# Update the selection to reflect the new extent of data entry
$ws.Range("A16").Select() | Out-Null
